# "update round 2 bidding" - refresh a handful of F (years_ihg) and
# H (allocated_number) cells on the `users` sheet with round-2 bidding
# results, then restore the workbook/sheet view to its post-edit state.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("users")

# --- Round 2 bidding data updates ---
# Column H (allocated_number) cells that used to hold the shared string
# "null" now hold the numeric allocation the user received in round 2.
# A few column F (years_ihg) cells were also corrected.
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(40, 8).Value = 9
$ws.Cells.Item(41, 8).Value = 5
$ws.Cells.Item(42, 8).Value = 6
$ws.Cells.Item(43, 8).Value = 14
$ws.Cells.Item(44, 6).Value = 6
$ws.Cells.Item(44, 8).Value = 1
$ws.Cells.Item(45, 8).Value = 24
$ws.Cells.Item(54, 8).Value = 95
$ws.Cells.Item(56, 8).Value = 28
$ws.Cells.Item(57, 8).Value = 27
$ws.Cells.Item(58, 8).Value = 16
$ws.Cells.Item(59, 8).Value = 17
$ws.Cells.Item(61, 8).Value = 45
$ws.Cells.Item(68, 8).Value = 88
$ws.Cells.Item(71, 6).Value = 2
$ws.Cells.Item(71, 8).Value = 17
$ws.Cells.Item(73, 8).Value = 30
$ws.Cells.Item(74, 8).Value = 6
$ws.Cells.Item(82, 8).Value = 16
$ws.Cells.Item(84, 8).Value = 42
$ws.Cells.Item(125, 8).Value = 4
$ws.Cells.Item(144, 8).Value = 12
$ws.Cells.Item(155, 8).Value = 3
$ws.Cells.Item(158, 8).Value = 8
$ws.Cells.Item(160, 8).Value = 98
$ws.Cells.Item(205, 8).Value = 20
$ws.Cells.Item(214, 6).Value = 2
$ws.Cells.Item(214, 8).Value = 41
$ws.Cells.Item(232, 8).Value = 18
$ws.Cells.Item(233, 8).Value = 15
$ws.Cells.Item(246, 8).Value = 55
$ws.Cells.Item(247, 8).Value = 20
$ws.Cells.Item(259, 6).Value = 6
$ws.Cells.Item(270, 8).Value = 8
$ws.Cells.Item(288, 8).Value = 15
$ws.Cells.Item(309, 8).Value = 71
$ws.Cells.Item(310, 8).Value = 14
$ws.Cells.Item(337, 8).Value = 31
$ws.Cells.Item(349, 8).Value = 19
$ws.Cells.Item(355, 8).Value = 22
$ws.Cells.Item(361, 8).Value = 37
$ws.Cells.Item(393, 8).Value = 22
$ws.Cells.Item(400, 8).Value = 96

# --- Column widths for E:H (years_ihg, points, has_mixed, allocated_number) ---
$ws.Range("E1:H1").ColumnWidth = 10

# --- View state: scroll back to top, freeze header row, select G19 ---
$ws.Range("A2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("G19").Select() | Out-Null
